# Vendor-SKU Map cleanup:
#  - Amazon "-FBA" duplicate SKU rows are removed (the Seller-Central and
#    Fulfilled-by-Amazon variants were being tracked as separate rows; the
#    FBA rows are no longer needed).
#  - The remaining Amazon "-SC" rows are renamed back to their base SKU
#    (the "SC" suffix is dropped) now that there is only one Amazon row
#    per SKU.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the surviving "-SC" suffixed SKUs (column B) back to their
#     base SKU name. Do this before any row deletion so the row numbers
#     below are unambiguous. ---
$ws.Cells.Item(307, 2).Value = "CordMate-OG"
$ws.Cells.Item(309, 2).Value = "CordMate-Pro"
$ws.Cells.Item(377, 2).Value = "SS100"
$ws.Cells.Item(378, 2).Value = "SD1_9"
$ws.Cells.Item(379, 2).Value = "SD8"
$ws.Cells.Item(380, 2).Value = "SD32"

# --- Delete the obsolete "-FBA" rows. Deleted from the bottom up so the
#     remaining row numbers referenced here stay valid while the sheet
#     shifts upward. ---
$ws.Rows(382).Delete()   # Amazon SD32FBA
$ws.Rows(381).Delete()   # Amazon SD8FBA
$ws.Rows(373).Delete()   # Amazon PPS24Kit-FBA
$ws.Rows(340).Delete()   # Amazon 6102FBA
$ws.Rows(310).Delete()   # Amazon CordMate-Pro-FBA
$ws.Rows(308).Delete()   # Amazon CordMate-OG-FBA

# --- Match the saved view state: plain selection on E308, no frozen /
#     scrolled topLeftCell override. ---
[void]$ws.Range("E308").Select()
